$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spreadsheet")

# --- Citation block (row 6): add the Title/DOI/Authors values under the headers in row 5 ---
# Set E10 first so the corrected "InchiCode" shared string is interned before the
# citation strings, matching the author's edit order.
$ws.Range("E10").Value = "InchiCode"

$ws.Range("D6").Value = "Examination of … "
$ws.Range("E6").Value = "DOI8444"
$ws.Range("F6").Value = "Author1, Author2"

# --- Clear the example Property/Value that had been pre-filled in the Datapoint grid ---
$ws.Range("D21").ClearContents()
$ws.Range("D23").ClearContents()

# --- Un-merge the title banner (A1:B2) and keep vertical centering only ---
$ws.Range("A1:B2").UnMerge()
$ws.Range("A1:B2").HorizontalAlignment = 1
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15.75

# --- Collapse the spacer rows between the Value row and the Variable section ---
$ws.Rows.Item(24).RowHeight = 1.5
$ws.Rows.Item(25).Hidden = $true
$ws.Rows.Item(26).RowHeight = 1.5
$ws.Rows.Item(27).RowHeight = 2.25
$ws.Rows.Item(28).RowHeight = 1.5

# --- Move the active selection to B2 ---
$ws.Activate()
$ws.Range("B2").Select()
